$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column Q (shifts Q:V -> S:X)
$insertRange = $ws.Range("Q1:R1").EntireColumn
$insertRange.Insert()

# Set the header values for the two newly inserted columns
$ws.Range("Q1").Value = "vali_c_matrix"
$ws.Range("R1").Value = "vali_c_matrix_perc"

# Rename the headers that moved / changed meaning
$ws.Range("V1").Value = "test_predicted"
$ws.Range("W1").Value = "test_c_matrix"
$ws.Range("X1").Value = "test_c_matrix_perc"
